# Fruta / hortaliza, semanal
# The weekly price-history rows (4-9) shift up: each row's data moves to the
# row two positions above it, wrapping around within the 4-9 block. In effect
# row 4 <- old row 6, row 5 <- old row 7, row 6 <- old row 8, row 7 <- old
# row 9, row 8 <- old row 4, row 9 <- old row 5. Only columns D, J, K, L, M,
# N, P, Q carry data that differs between rows; the rest (A, B, C, E, F, G,
# H, I, O, R) are already identical across the block and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRows = @(6, 7, 8, 9, 4, 5)
$targetRows = @(4, 5, 6, 7, 8, 9)

# Snapshot the "before" values for the columns that move, keyed by row.
# Use Value2 when reading back (Value returns a COM property descriptor
# in this host rather than the scalar) to get the raw number/string.
$snapshot = @{}
foreach ($r in @(4, 5, 6, 7, 8, 9)) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
    }
}

for ($i = 0; $i -lt $targetRows.Length; $i++) {
    $target = $targetRows[$i]
    $source = $sourceRows[$i]
    $data = $snapshot[$source]

    $ws.Cells.Item($target, 4).Value = $data.D
    $ws.Cells.Item($target, 10).Value = $data.J
    $ws.Cells.Item($target, 11).Value = $data.K
    $ws.Cells.Item($target, 12).Value = $data.L
    $ws.Cells.Item($target, 13).Value = $data.M
    $ws.Cells.Item($target, 14).Value = $data.N
    $ws.Cells.Item($target, 16).Value = $data.P
    $ws.Cells.Item($target, 17).Value = $data.Q
}
